$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New argument rows for the "add_internal_loads" measure (rows 43-74).
# Columns: A=Measure, B=Argument name, C=Argument display name, D=Type,
#          E=Unit, F=Description, G=Required, H=Default
$newRows = @(
    @("add_internal_loads", "electric_equipment_sched_weekday", "Electric equipment schedule for weekdays", "String", $null, "Schedule for the internal load of electric equipment for weekdays.", "True", $null),
    @("add_internal_loads", "electric_equipment_sched_saturday", "Electric equipment schedule for saturday", "String", $null, "Schedule for the internal load of electric equipment for holidays.", "True", $null),
    @("add_internal_loads", "electric_equipment_sched_sunday", "Electric equipment schedule for sunday", "String", $null, "Schedule for the internal load of electric equipment for sunday.", "True", $null),
    @("add_internal_loads", "electric_equipment_sched_holiday", "Electric equipment schedule for holidays", "String", $null, "Schedule for the internal load of electric equipment for holidays.", "True", $null),
    @("add_internal_loads", "lighting_sched_weekday", "Lighting schedule for weekdays", "String", $null, "Schedule for the internal load of lighting for weekdays.", "True", $null),
    @("add_internal_loads", "lighting_sched_saturday", "Lighting schedule for saturday", "String", $null, "Schedule for the internal load of lighting for saturday.", "True", $null),
    @("add_internal_loads", "lighting_sched_sunday", "Lighting schedule for sunday", "String", $null, "Schedule for the internal load of lighting for sunday.", "True", $null),
    @("add_internal_loads", "lighting_sched_holiday", "Lighting schedule for holiday", "String", $null, "Schedule for the internal load of lighting for holiday.", "True", $null),
    @("add_internal_loads", "people_sched_weekday", "People schedule for weekdays", "String", $null, "Schedule for the presence of people for weekdays.", "True", $null),
    @("add_internal_loads", "people_sched_saturday", "People schedule for saturday", "String", $null, "Schedule for the presence of people for saturday.", "True", $null),
    @("add_internal_loads", "people_sched_sunday", "People schedule for sunday", "String", $null, "Schedule for the presence of people for sunday.", "True", $null),
    @("add_internal_loads", "people_sched_holiday", "People schedule for holiday", "String", $null, "Schedule for the presence of people for holidays.", "True", $null),
    @("add_internal_loads", "people_activity_sched_weekday", "People activity schedule for weekdays", "String", $null, "Schedule for the activity of people for weekdays.", "True", $null),
    @("add_internal_loads", "people_activity_sched_satuday", "People activity schedule for saturday", "String", $null, "Schedule for the activity of people for saturday.", "True", $null),
    @("add_internal_loads", "people_activity_sched_sunday", "People activity schedule for sunday", "String", $null, "Schedule for the activity of people for sunday.", "True", $null),
    @("add_internal_loads", "people_activity_sched_holiday", "People activity schedule for holiday", "String", $null, "Schedule for the activity of people for holidays.", "True", $null),
    @("add_internal_loads", "nfa_gfa_ratio", "Ratio of NFA over GFA", "Float", $null, "Ratio of NFA over GFA.", "True", "1.0"),
    @("add_internal_loads", "electric_equipment_power_per_floor_area", "Area-specific electric equipment power", "Float", "kWh/m^-2", "Power of electric equipment relative to the GFA.", "True", $null),
    @("add_internal_loads", "lighting_power_per_floor_area", "Area-specific artificial lighting power", "Float", "kWh/m^-2", "Power of artificial lighting relative to the GFA.", "True", $null),
    @("add_internal_loads", "floor_area_per_person", "GFA per person", "Float", "m^2", "GFA per person.", "True", $null),
    @("add_internal_loads", "area_gfa_import", "GFA of imported model", "Float", "m^2", "(Export only) GFA of the imported model, if any.", "False", $null),
    @("add_internal_loads", "nfa_gfa_ratio_selection", "Selection of standard ratio of NFA over GFA", "String", $null, "(Export only) Selection of standard ratio of NFA over GFA.", "False", $null),
    @("add_internal_loads", "electric_equipment_sched_selection", "Selection of the electrip equipment schedule", "String", $null, "(Export only) Selection of the electric equipment schedule.", "False", $null),
    @("add_internal_loads", "lighting_sched_selection", "Selection of the artificial lighting schedule", "String", $null, "(Export only) Selection of the artificial lighting schedule.", "False", $null),
    @("add_internal_loads", "people_sched_selection", "Selection of the people schedule", "String", $null, "(Export only) Selection of the people presence schedule.", "False", $null),
    @("add_internal_loads", "people_activity_sched_selection", "Selection of the people activity schedule", "String", $null, "(Export only) Selection of the people activity schedule.", "False", $null),
    @("add_internal_loads", "is_custom_ratio", "Is custom NFA over GFA ratio", "Bool", $null, "(Export only) Flag whether the ratio of NFA over GFA is custom.", "False", $null),
    @("add_internal_loads", "is_imported_model", "Is imported model", "Bool", $null, "(Export only) Flag whether the geometric model is imported.", "False", $null),
    @("add_internal_loads", "is_custom_electric equipment", "Is custom electric equipment schedule", "Bool", $null, "(Export only) Flag whether the electric equipment schedule is custom.", "False", $null),
    @("add_internal_loads", "is_custom_lighting", "Is custom lighting schedule", "Bool", $null, "(Export only) Flag whether the artificial lighting schedule is custom.", "False", $null),
    @("add_internal_loads", "is_custom_people", "Is custom people schedule", "Bool", $null, "(Export only) Flag whether the people schedule is custom.", "False", $null),
    @("add_internal_loads", "is_custom_people_activity", "Is custom people activity", "Bool", $null, "(Export only) Flag whether the people activity schedule is custom.", "False", $null)
)

# Values that Excel's Range.Value setter would silently coerce into a
# native boolean/number (losing the "plain text" cell type). For these we
# write a text formula first, then convert the formula to a static value
# via Copy + PasteSpecial(xlPasteValues), which keeps the string type.
$textCoerced = @("True", "False", "1.0")

$startRow = 43
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 0; $c -lt 8; $c++) {
        $val = $row[$c]
        if ($null -eq $val) { continue }
        $cell = $ws.Cells.Item($r, $c + 1)
        if ($textCoerced -contains $val) {
            $cell.Formula = "=""" + $val + """"
            $cell.Copy() | Out-Null
            $cell.PasteSpecial(-4163) | Out-Null
        } else {
            $cell.Value = $val
        }
    }
}

$excel.CutCopyMode = 0

$ws.Range("B51").Select() | Out-Null
